# Schnittstelle_Server.xlsx - "Serverschnittstelle + localStorage-Nutzung aktualisiert."
#
# Adds a new API-table row documenting `game/abandon/{gameId}` (POST,
# returns 200/404, not yet implemented, requires SpringSecurity) to the
# "Tabelle4" ListObject on worksheet "Tabelle1", and renames the default
# cell style / theme to match the updated template (best effort - some of
# these are cosmetic workbook-level attributes that the host app may not
# expose through its object model, so we try them defensively).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- 1. Extend the "Tabelle4" table by one row (A4:K22 -> A4:K23) -----
$lo = $ws.ListObjects.Item("Tabelle4")
$lo.Resize($ws.Range("A4:K23"))

# --- 2. Copy the formatting of an existing, identically-styled data row
#        (row 10) down onto the freshly appended row 23 ----------------
$ws.Range("A10:K10").Copy()
$ws.Range("A23:K23").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 23 wraps onto two lines, same as the other multi-line rows.
$ws.Rows.Item(23).RowHeight = 30

# --- 3. Fill in the new row's values -----------------------------------
$ws.Range("A23").Value = "game/abandon/{gameId}"   # Resource
$ws.Range("C23").Value = "POST"                    # Method
$ws.Range("E23").Value = "200/404"                 # Returns / HTTP-Codes
$ws.Range("J23").Value = "nein"                     # Umgesetzt
$ws.Range("K23").Value = "ja"                       # SpringSecurity

# --- 4. The "Umgesetzt" dropdown validation now covers the new row -----
$ws.Range("J5:J23").Validation.Delete()
$ws.Range("J5:J23").Validation.Add(3, 1, 1, '"nein,ja,Mock"')

# --- 5. Keep the window selection/scroll position in sync with the
#        newly-added row ------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A24").Select()
try { $excel.ActiveWindow.ScrollRow = 19 } catch {}

# --- 6. Cosmetic rebrand: default cell style "Normal" -> "Standard" and
#        theme "Office" -> "Larissa" (best effort; not all hosts expose
#        these as writable object-model properties) ---------------------
try {
    $style = $wb.Styles.Item("Normal")
    $style.Name = "Standard"
} catch {}

try {
    $theme = $wb.Theme
    $theme.Name = "Larissa"
    try { $theme.ThemeColorScheme.Name = "Larissa" } catch {}
    try { $theme.ThemeFontScheme.Name = "Larissa" } catch {}
} catch {}
